$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New user rows appended to the master user_detail table (16th May Refresh)
$newRows = @(
    @{ Row = 34; Id = 110033; Uin = 9317596771; Name = "Nikola Tesla"; Email = "nikola.tesla@xyz.com"; Mobile = 818876434 },
    @{ Row = 35; Id = 110034; Uin = 9317596772; Name = "Graham Bell";  Email = "graham.bell@xyz.com";  Mobile = 818876435 },
    @{ Row = 36; Id = 110035; Uin = 9317596773; Name = "Albert Miles"; Email = "albert.miles@xyz.com"; Mobile = 818876436 }
)

# Duplicate the last existing data row (33) into each new row first, so the new
# rows inherit exactly the same cell formatting/styles used throughout the table
# (email column left-aligned-style, is_active column boolean style, etc.)
foreach ($r in $newRows) {
    $ws.Rows.Item(33).Copy()
    $ws.Rows.Item($r.Row).Insert(-4121)
}

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Id
    $ws.Cells.Item($row, 2).Value = $r.Uin
    $ws.Cells.Item($row, 5).Value = $r.Mobile
    $ws.Cells.Item($row, 6).Value = "ACT"
    $ws.Cells.Item($row, 7).Value = "eng"
    $ws.Cells.Item($row, 8).Value = "PWD"
    $ws.Cells.Item($row, 9).Value = $true
    $ws.Cells.Item($row, 10).Value = "superadmin"
    $ws.Cells.Item($row, 11).Value = "now()"
}

# Write the new (name, email) shared strings in the same order they appear in
# the target workbook: all three names first, then all three emails.
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Name
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.Email
}

# Move the active selection down to the row following the new data, matching the
# post-edit workbook's selection state
$ws.Range("A37:XFD1048576").Select()
